$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update title text (shared string) reflecting new reporting period (October -> November 2016)
$ws.Range("A2").Value = "Industrial Sector by Census Division and State, Year-to-Date through November 2016 (Continued)"

# Update relative standard error data values for rows 4-60 (columns E-I) per revised EPM run
    $ws.Range("H4").Value = 37
    $ws.Range("I4").Value = 14
    $ws.Range("I5").Value = 55
    $ws.Range("H6").Value = 37
    $ws.Range("I6").Value = 11
    $ws.Range("F7").Value = 173
    $ws.Range("I7").Value = 91
    $ws.Range("I8").Value = 209
    $ws.Range("E9").Value = 76
    $ws.Range("F9").Value = 6
    $ws.Range("I9").Value = 13
    $ws.Range("E10").Value = 191
    $ws.Range("F10").Value = 191
    $ws.Range("I10").Value = 46
    $ws.Range("E12").Value = 83
    $ws.Range("F12").Value = 9
    $ws.Range("I12").Value = 16
    $ws.Range("E13").Value = 138
    $ws.Range("F13").Value = 5
    $ws.Range("H13").Value = 11
    $ws.Range("H14").Value = 28
    $ws.Range("I14").Value = 13
    $ws.Range("F15").Value = 45
    $ws.Range("I15").Value = 12
    $ws.Range("F16").Value = 7
    $ws.Range("I16").Value = 16
    $ws.Range("E17").Value = 138
    $ws.Range("F17").Value = 13
    $ws.Range("I17").Value = 37
    $ws.Range("H18").Value = 66
    $ws.Range("I18").Value = 15
    $ws.Range("H19").Value = 58
    $ws.Range("I19").Value = 11
    $ws.Range("F20").Value = 46
    $ws.Range("H20").Value = 187
    $ws.Range("I20").Value = 13
    $ws.Range("F21").Value = 0
    $ws.Range("I21").Value = 112
    $ws.Range("H22").Value = 61
    $ws.Range("I22").Value = 18
    $ws.Range("F23").Value = 214
    $ws.Range("I23").Value = 153
    $ws.Range("I24").Value = 44
    $ws.Range("F25").Value = 0
    $ws.Range("I25").Value = 66
    $ws.Range("H26").Value = 6
    $ws.Range("I26").Value = 3
    $ws.Range("F27").Value = 99
    $ws.Range("F28").Value = 5
    $ws.Range("F29").Value = 2
    $ws.Range("I30").Value = 21
    $ws.Range("F31").Value = 4
    $ws.Range("H31").Value = 26
    $ws.Range("I31").Value = 10
    $ws.Range("I32").Value = 3
    $ws.Range("I33").Value = 8
    $ws.Range("I34").Value = 6
    $ws.Range("F35").Value = 2
    $ws.Range("H35").Value = 79
    $ws.Range("F36").Value = 3
    $ws.Range("I36").Value = 6
    $ws.Range("F37").Value = 2
    $ws.Range("I37").Value = 26
    $ws.Range("F38").Value = 2
    $ws.Range("H38").Value = 179
    $ws.Range("I38").Value = 9
    $ws.Range("F39").Value = 6
    $ws.Range("H40").Value = 12
    $ws.Range("I41").Value = 4
    $ws.Range("H42").Value = 15
    $ws.Range("F43").Value = 16
    $ws.Range("H43").Value = 102
    $ws.Range("I43").Value = 27
    $ws.Range("F44").Value = 8
    $ws.Range("H44").Value = 14
    $ws.Range("E45").Value = 169
    $ws.Range("F45").Value = 4
    $ws.Range("H45").Value = 15
    $ws.Range("F46").Value = 349
    $ws.Range("H46").Value = 70
    $ws.Range("I46").Value = 56
    $ws.Range("H47").Value = 49
    $ws.Range("I47").Value = 14
    $ws.Range("F48").Value = 54
    $ws.Range("I48").Value = 87
    $ws.Range("E49").Value = 169
    $ws.Range("F49").Value = 169
    $ws.Range("I49").Value = 20
    $ws.Range("I50").Value = 3020
    $ws.Range("I51").Value = 11
    $ws.Range("I52").Value = 12
    $ws.Range("E53").Value = 88
    $ws.Range("F53").Value = 4
    $ws.Range("H53").Value = 13
    $ws.Range("I53").Value = 2
    $ws.Range("E54").Value = 88
    $ws.Range("F54").Value = 11
    $ws.Range("H54").Value = 13
    $ws.Range("I54").Value = 3
    $ws.Range("F55").Value = 8
    $ws.Range("I55").Value = 10
    $ws.Range("F56").Value = 4
    $ws.Range("I56").Value = 4
    $ws.Range("F57").Value = 42
    $ws.Range("I57").Value = 37
    $ws.Range("F58").Value = 158
    $ws.Range("I58").Value = 66
    $ws.Range("F59").Value = 42
    $ws.Range("I59").Value = 44
    $ws.Range("E60").Value = 51
    $ws.Range("F60").Value = 1
    $ws.Range("H60").Value = 6
